$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-20 16:50:49"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-20 16:50:45"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-20 16:50:49"
